$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-lay the 4 columns: B=Email(new), C=NIF(was B), D=Colegio(was C); drop "Mesa" col ---
$nif2 = $ws.Range("B2").Value()
$nif3 = $ws.Range("B3").Value()
$nif4 = $ws.Range("B4").Value()
$colegio2 = $ws.Range("C2").Value()
$colegio3 = $ws.Range("C3").Value()
$colegio4 = $ws.Range("C4").Value()

$ws.Range("D1").Value = "Colegio"
$ws.Range("D2").Value = $colegio2
$ws.Range("D3").Value = $colegio3
$ws.Range("D4").Value = $colegio4

$ws.Range("C1").Value = "NIF"
$ws.Range("C2").Value = $nif2
$ws.Range("C3").Value = $nif3
$ws.Range("C4").Value = $nif4

$ws.Range("B1").Value = "Email"
$ws.Range("B2").Value = "jtp@hotmail.com"
$ws.Range("B3").Value = "llp@gmail.com"
$ws.Range("B4").Value = "atp@yahoo.com"

# --- Hyperlinks for the email column ---
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:jtp@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:llp@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:atp@yahoo.com")

$ws.Range("B5").Select()
